# Inserts 5 new mods (entity culling, ferritcore, immediatly fast,
# shulkerbox tooltip, skinlayers 3d) into the alphabetically sorted mods
# list on the active sheet, shifting the existing rows down as needed.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Shift existing rows' B:G content+format down to make room for
#        the 5 new rows. Processed from the bottom destination upward so
#        no source row is clobbered before it has been read. Column A
#        (the "verfuegbar" checkmark) is handled separately afterwards
#        because Copy() does not blank out a destination cell when the
#        source cell is empty, which would leave stale checkmarks behind.
#        Column I (the legend, rows 3:6 only) is left untouched. -------
$ws.Range("B29:G29").Copy($ws.Range("B34:G34"))
$ws.Range("B28:G28").Copy($ws.Range("B33:G33"))
$ws.Range("B27:G27").Copy($ws.Range("B32:G32"))
$ws.Range("B26:G26").Copy($ws.Range("B31:G31"))
$ws.Range("B25:G25").Copy($ws.Range("B30:G30"))
$ws.Range("B24:G24").Copy($ws.Range("B29:G29"))
$ws.Range("B23:G23").Copy($ws.Range("B28:G28"))
$ws.Range("B22:G22").Copy($ws.Range("B27:G27"))
$ws.Range("B21:G21").Copy($ws.Range("B26:G26"))
$ws.Range("B20:G20").Copy($ws.Range("B23:G23"))
$ws.Range("B19:G19").Copy($ws.Range("B22:G22"))
$ws.Range("B18:G18").Copy($ws.Range("B21:G21"))
$ws.Range("B17:G17").Copy($ws.Range("B20:G20"))
$ws.Range("B16:G16").Copy($ws.Range("B19:G19"))
$ws.Range("B15:G15").Copy($ws.Range("B18:G18"))
$ws.Range("B14:G14").Copy($ws.Range("B17:G17"))
$ws.Range("B13:G13").Copy($ws.Range("B14:G14"))
$ws.Range("B12:G12").Copy($ws.Range("B13:G13"))

# --- 2. The 3 brand-new rows (12, 24, 25) need their B:G formatting
#        reset to the plain/default look. Row 12 and row 24 already
#        carry the default C:G styling (s=7) from the stale data that
#        used to live there, but row 25 inherited the special
#        "worldeditcui" E/G styling (s=10 / s=11) which must be reset to
#        the plain style (copied from the always-plain E7:G7). ---------
$ws.Range("E7:G7").Copy($ws.Range("E25:G25"))

# --- 3. Write the new mod names into column B. The order below matches
#        the order the strings were appended to the shared-string table
#        so the resulting indices line up with the target file. -------
$ws.Range("B12").Value = "entity culling"
$ws.Range("B16").Value = "immediatly fast"
$ws.Range("B15").Value = "ferritcore"
$ws.Range("B24").Value = "shulkerbox tooltip"
$ws.Range("B25").Value = "skinlayers 3d"

# --- 4. Column A "verfuegbar" checkmarks for rows 12:34, set explicitly
#        (rather than copied) for the reason noted in step 1. Rows 30:34
#        did not exist before this edit, so column A there first needs
#        the plain cell style (copied from A7) before a value is set,
#        otherwise the new cells come out unstyled. -------------------
$ws.Range("A7").Copy($ws.Range("A12:A34"))
$ws.Range("A12").ClearContents()
$ws.Range("A13").Value = "x"
$ws.Range("A14").Value = "x"
$ws.Range("A15").Value = "x"
$ws.Range("A16").Value = "x"
$ws.Range("A17").Value = "x"
$ws.Range("A18").ClearContents()
$ws.Range("A19").Value = "x"
$ws.Range("A20").Value = "x"
$ws.Range("A21").Value = "x"
$ws.Range("A22").ClearContents()
$ws.Range("A23").ClearContents()
$ws.Range("A24").ClearContents()
$ws.Range("A25").ClearContents()
$ws.Range("A26").Value = "x"
$ws.Range("A27").Value = "x"
$ws.Range("A28").Value = "x"
$ws.Range("A29").ClearContents()
$ws.Range("A30").ClearContents()
$ws.Range("A31").Value = "x"
$ws.Range("A32").Value = "x"
$ws.Range("A33").Value = "x"
$ws.Range("A34").Value = "x"

# --- 5. Update the sheet's scroll position / selection ----------------
$ws.Application.ActiveWindow.ScrollRow = 4
$ws.Range("B25").Select()
